$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "65.594.05"
Set-TextValue "E2" "  +2.64%  "

Set-TextValue "D3" "3.189.86"
Set-TextValue "E3" "  +4.50%  "

Set-TextValue "E4" "  +0.05%  "

Set-TextValue "D5" "574.81"
Set-TextValue "E5" "  +3.17%  "

Set-TextValue "D6" "152.76"
Set-TextValue "E6" "  +7.45%  "

Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  +0.01%  "

Set-TextValue "D8" "3.185.92"
Set-TextValue "E8" "  +4.44%  "

Set-TextValue "D9" "0.530"
Set-TextValue "E9" "  +4.60%  "

Set-TextValue "D10" "0.165"
Set-TextValue "E10" "  +6.32%  "

Set-TextValue "D11" "6.23"
Set-TextValue "E11" "  +2.84%  "

Set-TextValue "D12" "0.511"
Set-TextValue "E12" "  +7.34%  "

Set-TextValue "D13" "0.0000276"
Set-TextValue "E13" "  +19.93%  "

Set-TextValue "D14" "38.61"
Set-TextValue "E14" "  +10.65%  "

Set-TextValue "D15" "3.711.10"
Set-TextValue "E15" "  +4.51%  "

Set-TextValue "D16" "65.588.29"
Set-TextValue "E16" "  +2.60%  "

Set-TextValue "D17" "7.27"
Set-TextValue "E17" "  +8.43%  "

Set-TextValue "D18" "3.188.96"
Set-TextValue "E18" "  +4.37%  "

Set-TextValue "E19" "  +1.40%  "

Set-TextValue "D20" "517.18"
Set-TextValue "E20" "  +9.02%  "

Set-TextValue "D21" "15.02"
Set-TextValue "E21" "  +7.26%  "

Set-TextValue "D22" "16.03"
Set-TextValue "E22" "  +12.86%  "

Set-TextValue "D23" "0.743"
Set-TextValue "E23" "  +9.90%  "

Set-TextValue "D24" "7.93"
Set-TextValue "E24" "  +5.15%  "

Set-TextValue "D25" "85.22"
Set-TextValue "E25" "  +4.70%  "

Set-TextValue "E26" "  +0.17%  "

Set-TextValue "D27" "9.17"
Set-TextValue "E27" "  +16.10%  "

Set-TextValue "D28" "2.93"
Set-TextValue "E28" "  +5.15%  "

Set-TextValue "E29" "  +10.10%  "

Set-TextValue "D30" "28.27"
Set-TextValue "E30" "  +7.85%  "

Set-TextValue "E31" "  +16.79%  "

Set-TextValue "D32" "1.24"
Set-TextValue "E32" "  +8.45%  "

Set-TextValue "D33" "1.00"
Set-TextValue "E33" "  -0.04%  "

Set-TextValue "D34" "6.32"
Set-TextValue "E34" "  +13.51%  "

Set-TextValue "D35" "6.68"
Set-TextValue "E35" "  +7.99%  "

Set-TextValue "D36" "55.85"
Set-TextValue "E36" "  +1.96%  "

Set-TextValue "D37" "485.43"
Set-TextValue "E37" "  +9.86%  "

Set-TextValue "B38" "dogwifhat"
Set-TextValue "C38" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D38" "3.19"
Set-TextValue "E38" "  +13.15%  "

Set-TextValue "B39" "Hedera"
Set-TextValue "C39" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D39" "0.0883"
Set-TextValue "E39" "  +9.84%  "

Set-TextValue "D40" "0.0425"
Set-TextValue "E40" "  +5.29%  "

Set-TextValue "D41" "3.152.86"
Set-TextValue "E41" "  +6.48%  "

Set-TextValue "D42" "8.70"
Set-TextValue "E42" "  +6.29%  "

Set-TextValue "E43" "  +7.51%  "

Set-TextValue "D44" "2.52"
Set-TextValue "E44" "  +17.83%  "

Set-TextValue "E45" "  +12.56%  "

Set-TextValue "D46" "29.53"
Set-TextValue "E46" "  +6.97%  "

Set-TextValue "D47" "0.0₃0591"
Set-TextValue "E47" "  +15.65%  "

Set-TextValue "E48" "  -0.05%  "

Set-TextValue "E49" "  +3.57%  "

Set-TextValue "E50" "  +14.01%  "

Set-TextValue "D51" "124.07"
Set-TextValue "E51" "  +5.92%  "
